$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.772158333333334
$ws.Range("H2").Value = 14.316475
$ws.Range("I2").Value = 0.2444103987677044
$ws.Range("J2").Value = 0.2444103987677044
$ws.Range("M2").Value = 51.956163
$ws.Range("N2").Value = 155.868489
$ws.Range("O2").Value = 0.2075159764120945
$ws.Range("P2").Value = 0.2075159764120945
$ws.Range("Q2").Value = 247.943036228475
$ws.Range("R2").Value = 2231.487326056275
$ws.Range("S2").Value = 0.05071906254554956
$ws.Range("T2").Value = 0.05071906254554956
# Row 3
$ws.Range("G3").Value = 4.772158333333334
$ws.Range("H3").Value = 14.316475
$ws.Range("I3").Value = 0.2444103987677044
$ws.Range("J3").Value = 0.2444103987677044
$ws.Range("O3").Value = 0.5961756600906958
$ws.Range("P3").Value = 0.5961756600906958
$ws.Range("Q3").Value = 712.3191469116556
$ws.Range("R3").Value = 6410.8723222049
$ws.Range("S3").Value = 0.1457115308183664
$ws.Range("T3").Value = 0.1457115308183664
# Row 4
$ws.Range("G4").Value = 4.772158333333334
$ws.Range("H4").Value = 14.316475
$ws.Range("I4").Value = 0.2444103987677044
$ws.Range("J4").Value = 0.2444103987677044
$ws.Range("M4").Value = 15.539306
$ws.Range("N4").Value = 46.617918
$ws.Range("O4").Value = 0.06206490377967901
$ws.Range("P4").Value = 0.06206490377967901
$ws.Range("Q4").Value = 74.15602862211668
$ws.Range("R4").Value = 667.4042575990501
$ws.Range("S4").Value = 0.01516930788227055
$ws.Range("T4").Value = 0.01516930788227055
# Row 5
$ws.Range("G5").Value = 4.772158333333334
$ws.Range("H5").Value = 14.316475
$ws.Range("I5").Value = 0.2444103987677044
$ws.Range("J5").Value = 0.2444103987677044
$ws.Range("M5").Value = 15.68808233333333
$ws.Range("N5").Value = 47.064247
$ws.Range("O5").Value = 0.06265912522129466
$ws.Range("P5").Value = 0.06265912522129466
$ws.Range("Q5").Value = 74.86601284103611
$ws.Range("R5").Value = 673.7941155693251
$ws.Range("S5").Value = 0.01531454178177215
$ws.Range("T5").Value = 0.01531454178177215
# Row 6
$ws.Range("G6").Value = 4.772158333333334
$ws.Range("H6").Value = 14.316475
$ws.Range("I6").Value = 0.2444103987677044
$ws.Range("J6").Value = 0.2444103987677044
$ws.Range("M6").Value = 3.763510333333333
$ws.Range("N6").Value = 11.290531
$ws.Range("O6").Value = 0.01503168202699406
$ws.Range("P6").Value = 0.01503168202699406
$ws.Range("Q6").Value = 17.96006719980278
$ws.Range("R6").Value = 161.640604798225
$ws.Range("S6").Value = 0.003673899398366953
$ws.Range("T6").Value = 0.003673899398366953
# Row 7
$ws.Range("G7").Value = 4.772158333333334
$ws.Range("H7").Value = 14.316475
$ws.Range("I7").Value = 0.2444103987677044
$ws.Range("J7").Value = 0.2444103987677044
$ws.Range("M7").Value = 14.15919333333333
$ws.Range("N7").Value = 42.47758
$ws.Range("O7").Value = 0.05655265246924192
$ws.Range("P7").Value = 0.05655265246924192
$ws.Range("Q7").Value = 67.56991245894446
$ws.Range("R7").Value = 608.1292121305
$ws.Range("S7").Value = 0.01382205634137882
$ws.Range("T7").Value = 0.01382205634137882
# Row 8
$ws.Range("G8").Value = 11.54138666666667
$ws.Range("H8").Value = 34.62416
$ws.Range("I8").Value = 0.591102541135077
$ws.Range("J8").Value = 0.591102541135077
$ws.Range("M8").Value = 51.956163
$ws.Range("N8").Value = 155.868489
$ws.Range("O8").Value = 0.2075159764120945
$ws.Range("P8").Value = 0.2075159764120945
$ws.Range("Q8").Value = 599.6461668993601
$ws.Range("R8").Value = 5396.815502094241
$ws.Range("S8").Value = 0.1226632209833158
$ws.Range("T8").Value = 0.1226632209833158
# Row 9
$ws.Range("G9").Value = 11.54138666666667
$ws.Range("H9").Value = 34.62416
$ws.Range("I9").Value = 0.591102541135077
$ws.Range("J9").Value = 0.591102541135077
$ws.Range("O9").Value = 0.5961756600906958
$ws.Range("P9").Value = 0.5961756600906958
$ws.Range("Q9").Value = 1722.732174905671
$ws.Range("R9").Value = 15504.58957415104
$ws.Range("S9").Value = 0.3524009476424922
$ws.Range("T9").Value = 0.3524009476424922
# Row 10
$ws.Range("G10").Value = 11.54138666666667
$ws.Range("H10").Value = 34.62416
$ws.Range("I10").Value = 0.591102541135077
$ws.Range("J10").Value = 0.591102541135077
$ws.Range("M10").Value = 15.539306
$ws.Range("N10").Value = 46.617918
$ws.Range("O10").Value = 0.06206490377967901
$ws.Range("P10").Value = 0.06206490377967901
$ws.Range("Q10").Value = 179.3451390776534
$ws.Range("R10").Value = 1614.10625169888
$ws.Range("S10").Value = 0.03668672233947231
$ws.Range("T10").Value = 0.0366867223394723
# Row 11
$ws.Range("G11").Value = 11.54138666666667
$ws.Range("H11").Value = 34.62416
$ws.Range("I11").Value = 0.591102541135077
$ws.Range("J11").Value = 0.591102541135077
$ws.Range("M11").Value = 15.68808233333333
$ws.Range("N11").Value = 47.064247
$ws.Range("O11").Value = 0.06265912522129466
$ws.Range("P11").Value = 0.06265912522129466
$ws.Range("Q11").Value = 181.0622242675022
$ws.Range("R11").Value = 1629.56001840752
$ws.Range("S11").Value = 0.03703796814360827
$ws.Range("T11").Value = 0.03703796814360827
# Row 12
$ws.Range("G12").Value = 11.54138666666667
$ws.Range("H12").Value = 34.62416
$ws.Range("I12").Value = 0.591102541135077
$ws.Range("J12").Value = 0.591102541135077
$ws.Range("M12").Value = 3.763510333333333
$ws.Range("N12").Value = 11.290531
$ws.Range("O12").Value = 0.01503168202699406
$ws.Range("P12").Value = 0.01503168202699406
$ws.Range("Q12").Value = 43.43612798099555
$ws.Range("R12").Value = 390.9251518289599
$ws.Range("S12").Value = 0.008885265443690653
$ws.Range("T12").Value = 0.008885265443690653
# Row 13
$ws.Range("G13").Value = 11.54138666666667
$ws.Range("H13").Value = 34.62416
$ws.Range("I13").Value = 0.591102541135077
$ws.Range("J13").Value = 0.591102541135077
$ws.Range("M13").Value = 14.15919333333333
$ws.Range("N13").Value = 42.47758
$ws.Range("O13").Value = 0.05655265246924192
$ws.Range("P13").Value = 0.05655265246924192
$ws.Range("Q13").Value = 163.4167251480889
$ws.Range("R13").Value = 1470.7505263328
$ws.Range("S13").Value = 0.03342841658249779
$ws.Range("T13").Value = 0.03342841658249779
# Row 14
$ws.Range("G14").Value = 3.211640333333333
$ws.Range("H14").Value = 9.634920999999999
$ws.Range("I14").Value = 0.1644870600972187
$ws.Range("J14").Value = 0.1644870600972187
$ws.Range("M14").Value = 51.956163
$ws.Range("N14").Value = 155.868489
$ws.Range("O14").Value = 0.2075159764120945
$ws.Range("P14").Value = 0.2075159764120945
$ws.Range("Q14").Value = 166.864508656041
$ws.Range("R14").Value = 1501.780577904369
$ws.Range("S14").Value = 0.03413369288322921
$ws.Range("T14").Value = 0.03413369288322921
# Row 15
$ws.Range("G15").Value = 3.211640333333333
$ws.Range("H15").Value = 9.634920999999999
$ws.Range("I15").Value = 0.1644870600972187
$ws.Range("J15").Value = 0.1644870600972187
$ws.Range("O15").Value = 0.5961756600906958
$ws.Range("P15").Value = 0.5961756600906958
$ws.Range("Q15").Value = 479.3874684432582
$ws.Range("R15").Value = 4314.487215989324
$ws.Range("S15").Value = 0.09806318162983731
$ws.Range("T15").Value = 0.09806318162983731
# Row 16
$ws.Range("G16").Value = 3.211640333333333
$ws.Range("H16").Value = 9.634920999999999
$ws.Range("I16").Value = 0.1644870600972187
$ws.Range("J16").Value = 0.1644870600972187
$ws.Range("M16").Value = 15.539306
$ws.Range("N16").Value = 46.617918
$ws.Range("O16").Value = 0.06206490377967901
$ws.Range("P16").Value = 0.06206490377967901
$ws.Range("Q16").Value = 49.90666190160866
$ws.Range("R16").Value = 449.1599571144779
$ws.Range("S16").Value = 0.01020887355793616
$ws.Range("T16").Value = 0.01020887355793616
# Row 17
$ws.Range("G17").Value = 3.211640333333333
$ws.Range("H17").Value = 9.634920999999999
$ws.Range("I17").Value = 0.1644870600972187
$ws.Range("J17").Value = 0.1644870600972187
$ws.Range("M17").Value = 15.68808233333333
$ws.Range("N17").Value = 47.064247
$ws.Range("O17").Value = 0.06265912522129466
$ws.Range("P17").Value = 0.06265912522129466
$ws.Range("Q17").Value = 50.38447797438744
$ws.Range("R17").Value = 453.4603017694869
$ws.Range("S17").Value = 0.01030661529591425
$ws.Range("T17").Value = 0.01030661529591425
# Row 18
$ws.Range("G18").Value = 3.211640333333333
$ws.Range("H18").Value = 9.634920999999999
$ws.Range("I18").Value = 0.1644870600972187
$ws.Range("J18").Value = 0.1644870600972187
$ws.Range("M18").Value = 3.763510333333333
$ws.Range("N18").Value = 11.290531
$ws.Range("O18").Value = 0.01503168202699406
$ws.Range("P18").Value = 0.01503168202699406
$ws.Range("Q18").Value = 12.08704158145011
$ws.Range("R18").Value = 108.783374233051
$ws.Range("S18").Value = 0.002472517184936454
$ws.Range("T18").Value = 0.002472517184936454
# Row 19
$ws.Range("G19").Value = 3.211640333333333
$ws.Range("H19").Value = 9.634920999999999
$ws.Range("I19").Value = 0.1644870600972187
$ws.Range("J19").Value = 0.1644870600972187
$ws.Range("M19").Value = 14.15919333333333
$ws.Range("N19").Value = 42.47758
$ws.Range("O19").Value = 0.05655265246924192
$ws.Range("P19").Value = 0.05655265246924192
$ws.Range("Q19").Value = 45.47423639679778
$ws.Range("R19").Value = 409.2681275711799
$ws.Range("S19").Value = 0.00930217954536532
$ws.Range("T19").Value = 0.00930217954536532
